$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove existing hyperlinks before rewriting data to avoid stale links
$ws.Hyperlinks.Delete()

# Clear old data rows (header stays in row 1)
$ws.Range("A2:H23").ClearContents()

# Write refreshed data set (27 listings), sorted by score desc, new scrape timestamp
$ws.Cells.Item(2,1).Value = '2025-11-12 18:25:55'
$ws.Cells.Item(2,2).Value = '専門データ分析:AIコスト最適化設計と厳格な機密保持を必須とするWebシステム開発(段階的継続発注)'
$ws.Cells.Item(2,3).Value = 'システム開発'
$ws.Cells.Item(2,4).Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Cells.Item(2,5).Value = '期限情報なし'
$ws.Cells.Item(2,6).Value = 'https://www.lancers.jp/work/detail/5431917'
$ws.Cells.Item(2,7).Value = 403
$ws.Cells.Item(2,8).Value = '🔥AI,Ai ◆開発,システム開発'

$ws.Cells.Item(3,1).Value = '2025-11-12 18:25:55'
$ws.Cells.Item(3,2).Value = '詳細設計及び、Next.js,node.jsによるWEBアプリケーション開発'
$ws.Cells.Item(3,3).Value = 'システム開発'
$ws.Cells.Item(3,4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(3,5).Value = '期限情報なし'
$ws.Cells.Item(3,6).Value = 'https://www.lancers.jp/work/detail/5427010'
$ws.Cells.Item(3,7).Value = 245
$ws.Cells.Item(3,8).Value = '🔥Next.js ◆開発,Node.js ◇アプリ'

$ws.Cells.Item(4,1).Value = '2025-11-12 18:25:55'
$ws.Cells.Item(4,2).Value = '<Next.js、バックエンド開発> ガントチャートアプリの改修製造'
$ws.Cells.Item(4,3).Value = 'システム開発'
$ws.Cells.Item(4,4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(4,5).Value = '期限情報なし'
$ws.Cells.Item(4,6).Value = 'https://www.lancers.jp/work/detail/5427011'
$ws.Cells.Item(4,7).Value = 225
$ws.Cells.Item(4,8).Value = '🔥Next.js ◆開発 ◇アプリ'

$ws.Cells.Item(5,1).Value = '2025-11-12 18:25:55'
$ws.Cells.Item(5,2).Value = '【急募】大手保険会社向けスマホアプリ設計書作成依頼'
$ws.Cells.Item(5,3).Value = 'システム開発'
$ws.Cells.Item(5,4).Value = '100,000 円 ~ 200,000 円 / 募集期間 2 日、取引期間 0 日'
$ws.Cells.Item(5,5).Value = '期限情報なし'
$ws.Cells.Item(5,6).Value = 'https://www.lancers.jp/work/detail/5431609'
$ws.Cells.Item(5,7).Value = 103
$ws.Cells.Item(5,8).Value = '★スマホアプリ ◇アプリ'

$ws.Cells.Item(6,1).Value = '2025-11-12 18:25:55'
$ws.Cells.Item(6,2).Value = 'WEBサイトへの自動ログインのツール'
$ws.Cells.Item(6,3).Value = 'システム開発'
$ws.Cells.Item(6,4).Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Cells.Item(6,5).Value = '期限情報なし'
$ws.Cells.Item(6,6).Value = 'https://www.lancers.jp/work/detail/5432620'
$ws.Cells.Item(6,7).Value = 90
$ws.Cells.Item(6,8).Value = '◆ツール ◇サイト'

$ws.Cells.Item(7,1).Value = '2025-11-12 18:25:55'
$ws.Cells.Item(7,2).Value = '【アセスメント試験】Microsoft Formsデータ処理の自動化依頼'
$ws.Cells.Item(7,3).Value = 'システム開発'
$ws.Cells.Item(7,4).Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Cells.Item(7,5).Value = '期限情報なし'
$ws.Cells.Item(7,6).Value = 'https://www.lancers.jp/work/detail/5432776'
$ws.Cells.Item(7,7).Value = 88
$ws.Cells.Item(7,8).Value = '◆自動化'

$ws.Cells.Item(8,1).Value = '2025-11-12 18:25:55'
$ws.Cells.Item(8,2).Value = '【日本人限定・長期募集】SNS運用担当募集|Web開発会社 JapanDream'
$ws.Cells.Item(8,3).Value = 'システム開発'
$ws.Cells.Item(8,4).Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Cells.Item(8,5).Value = '期限情報なし'
$ws.Cells.Item(8,6).Value = 'https://www.lancers.jp/work/detail/5432819'
$ws.Cells.Item(8,7).Value = 68
$ws.Cells.Item(8,8).Value = '◆開発'

$ws.Cells.Item(9,1).Value = '2025-11-12 18:25:55'
$ws.Cells.Item(9,2).Value = 'ヤフオクの指定出品者の出品物を一括してウォッチリストに登録するツール'
$ws.Cells.Item(9,3).Value = 'システム開発'
$ws.Cells.Item(9,4).Value = '~ 5,000 円 / 固定'
$ws.Cells.Item(9,5).Value = '期限情報なし'
$ws.Cells.Item(9,6).Value = 'https://www.lancers.jp/work/detail/5431786'
$ws.Cells.Item(9,7).Value = 65
$ws.Cells.Item(9,8).Value = '◆ツール'

$ws.Cells.Item(10,1).Value = '2025-11-12 18:25:55'
$ws.Cells.Item(10,2).Value = 'PHP業務アプリケーションの改修対応'
$ws.Cells.Item(10,3).Value = 'システム開発'
$ws.Cells.Item(10,4).Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Cells.Item(10,5).Value = '期限情報なし'
$ws.Cells.Item(10,6).Value = 'https://www.lancers.jp/work/detail/5426598'
$ws.Cells.Item(10,7).Value = 58
$ws.Cells.Item(10,8).Value = '◇アプリ ○PHP'

$ws.Cells.Item(11,1).Value = '2025-11-12 18:25:55'
$ws.Cells.Item(11,2).Value = '【案件】既存WordPressサイトの読み込み速度改善'
$ws.Cells.Item(11,3).Value = 'システム開発'
$ws.Cells.Item(11,4).Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Cells.Item(11,5).Value = '期限情報なし'
$ws.Cells.Item(11,6).Value = 'https://www.lancers.jp/work/detail/5432161'
$ws.Cells.Item(11,7).Value = 50
$ws.Cells.Item(11,8).Value = '◇サイト ○WordPress'

$ws.Cells.Item(12,1).Value = '2025-11-12 18:25:55'
$ws.Cells.Item(12,2).Value = '【急募】ショッピファイでジャケット仕様確定システム構築'
$ws.Cells.Item(12,3).Value = 'システム開発'
$ws.Cells.Item(12,4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(12,5).Value = '期限情報なし'
$ws.Cells.Item(12,6).Value = 'https://www.lancers.jp/work/detail/5432465'
$ws.Cells.Item(12,7).Value = 40

$ws.Cells.Item(13,1).Value = '2025-11-12 18:25:55'
$ws.Cells.Item(13,2).Value = '小売店向けシステム性能試験'
$ws.Cells.Item(13,3).Value = 'システム開発'
$ws.Cells.Item(13,4).Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Cells.Item(13,5).Value = '期限情報なし'
$ws.Cells.Item(13,6).Value = 'https://www.lancers.jp/work/detail/5430176'
$ws.Cells.Item(13,7).Value = 40

$ws.Cells.Item(14,1).Value = '2025-11-12 18:25:55'
$ws.Cells.Item(14,2).Value = '【電気錠制御】オフィスセキュリティシステム刷新の協力者募集'
$ws.Cells.Item(14,3).Value = 'システム開発'
$ws.Cells.Item(14,4).Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Cells.Item(14,5).Value = '期限情報なし'
$ws.Cells.Item(14,6).Value = 'https://www.lancers.jp/work/detail/5431852'
$ws.Cells.Item(14,7).Value = 33

$ws.Cells.Item(15,1).Value = '2025-11-12 18:25:55'
$ws.Cells.Item(15,2).Value = 'Networkエンジニア'
$ws.Cells.Item(15,3).Value = 'システム開発'
$ws.Cells.Item(15,4).Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Cells.Item(15,5).Value = '期限情報なし'
$ws.Cells.Item(15,6).Value = 'https://www.lancers.jp/work/detail/5432661'
$ws.Cells.Item(15,7).Value = 25

$ws.Cells.Item(16,1).Value = '2025-11-12 18:25:55'
$ws.Cells.Item(16,2).Value = '【急募】神奈川県藤沢市でNWエンジニアを募集!(2.5ヶ月)'
$ws.Cells.Item(16,3).Value = 'システム開発'
$ws.Cells.Item(16,4).Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Cells.Item(16,5).Value = '期限情報なし'
$ws.Cells.Item(16,6).Value = 'https://www.lancers.jp/work/detail/5432622'
$ws.Cells.Item(16,7).Value = 25

$ws.Cells.Item(17,1).Value = '2025-11-12 18:25:55'
$ws.Cells.Item(17,2).Value = '【フルリモート/継続案件】Salesforceの導入・運用・保守'
$ws.Cells.Item(17,3).Value = 'システム開発'
$ws.Cells.Item(17,4).Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Cells.Item(17,5).Value = '期限情報なし'
$ws.Cells.Item(17,6).Value = 'https://www.lancers.jp/work/detail/5432563'
$ws.Cells.Item(17,7).Value = 25

$ws.Cells.Item(18,1).Value = '2025-11-12 18:25:55'
$ws.Cells.Item(18,2).Value = 'OR(operations research)にて最適化の仕組みの構築 (リモート)'
$ws.Cells.Item(18,3).Value = 'システム開発'
$ws.Cells.Item(18,4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(18,5).Value = '期限情報なし'
$ws.Cells.Item(18,6).Value = 'https://www.lancers.jp/work/detail/5427007'
$ws.Cells.Item(18,7).Value = 25

$ws.Cells.Item(19,1).Value = '2025-11-12 18:25:55'
$ws.Cells.Item(19,2).Value = 'OR(operations research)にて最適化の仕組みの構築(社内常駐)'
$ws.Cells.Item(19,3).Value = 'システム開発'
$ws.Cells.Item(19,4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(19,5).Value = '期限情報なし'
$ws.Cells.Item(19,6).Value = 'https://www.lancers.jp/work/detail/5427009'
$ws.Cells.Item(19,7).Value = 25

$ws.Cells.Item(20,1).Value = '2025-11-12 18:25:55'
$ws.Cells.Item(20,2).Value = '【急募】楽天市場在庫連動システム(同一店舗内)のエラー修正依頼'
$ws.Cells.Item(20,3).Value = 'システム開発'
$ws.Cells.Item(20,4).Value = '10,000 円 ~ 20,000 円 / 固定'
$ws.Cells.Item(20,5).Value = '期限情報なし'
$ws.Cells.Item(20,6).Value = 'https://www.lancers.jp/work/detail/5432212'
$ws.Cells.Item(20,7).Value = 25

$ws.Cells.Item(21,1).Value = '2025-11-12 18:25:55'
$ws.Cells.Item(21,2).Value = '適合商品検索ページ作成'
$ws.Cells.Item(21,3).Value = 'システム開発'
$ws.Cells.Item(21,4).Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Cells.Item(21,5).Value = '期限情報なし'
$ws.Cells.Item(21,6).Value = 'https://www.lancers.jp/work/detail/5432621'
$ws.Cells.Item(21,7).Value = 18

$ws.Cells.Item(22,1).Value = '2025-11-12 18:25:55'
$ws.Cells.Item(22,2).Value = '【急募】企業向け情シス・セキュリティ業務支援'
$ws.Cells.Item(22,3).Value = 'システム開発'
$ws.Cells.Item(22,4).Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Cells.Item(22,5).Value = '期限情報なし'
$ws.Cells.Item(22,6).Value = 'https://www.lancers.jp/work/detail/5432712'
$ws.Cells.Item(22,7).Value = 18

$ws.Cells.Item(23,1).Value = '2025-11-12 18:25:55'
$ws.Cells.Item(23,2).Value = '【音楽制作】サイケデリックトランスのトラックを作成してくれる方募集'
$ws.Cells.Item(23,3).Value = 'システム開発'
$ws.Cells.Item(23,4).Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Cells.Item(23,5).Value = '期限情報なし'
$ws.Cells.Item(23,6).Value = 'https://www.lancers.jp/work/detail/5432042'
$ws.Cells.Item(23,7).Value = 18

$ws.Cells.Item(24,1).Value = '2025-11-12 18:25:55'
$ws.Cells.Item(24,2).Value = '【スポット案件】HTML途切れ・白画面・Segmentation fault調査対応'
$ws.Cells.Item(24,3).Value = 'システム開発'
$ws.Cells.Item(24,4).Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Cells.Item(24,5).Value = '期限情報なし'
$ws.Cells.Item(24,6).Value = 'https://www.lancers.jp/work/detail/5432323'
$ws.Cells.Item(24,7).Value = 13

$ws.Cells.Item(25,1).Value = '2025-11-12 18:25:55'
$ws.Cells.Item(25,2).Value = '【急募】Wartalesの武器アイコンとモデルを日本刀に差し替え'
$ws.Cells.Item(25,3).Value = 'システム開発'
$ws.Cells.Item(25,4).Value = '10,000 円 ~ 20,000 円 / 固定'
$ws.Cells.Item(25,5).Value = '期限情報なし'
$ws.Cells.Item(25,6).Value = 'https://www.lancers.jp/work/detail/5432425'
$ws.Cells.Item(25,7).Value = 10

$ws.Cells.Item(26,1).Value = '2025-11-12 18:25:55'
$ws.Cells.Item(26,2).Value = '初回 MT4用インジケータの修正カスタマイズ(.mq4)'
$ws.Cells.Item(26,3).Value = 'システム開発'
$ws.Cells.Item(26,4).Value = '10,000 円 ~ 20,000 円 / 固定'
$ws.Cells.Item(26,5).Value = '期限情報なし'
$ws.Cells.Item(26,6).Value = 'https://www.lancers.jp/work/detail/5432362'
$ws.Cells.Item(26,7).Value = 10

$ws.Cells.Item(27,1).Value = '2025-11-12 18:25:55'
$ws.Cells.Item(27,2).Value = 'MT4用インジケータの修正カスタマイズ(.mq4)'
$ws.Cells.Item(27,3).Value = 'システム開発'
$ws.Cells.Item(27,4).Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Cells.Item(27,5).Value = '期限情報なし'
$ws.Cells.Item(27,6).Value = 'https://www.lancers.jp/work/detail/5432305'
$ws.Cells.Item(27,7).Value = 10

$ws.Cells.Item(28,1).Value = '2025-11-12 18:25:55'
$ws.Cells.Item(28,2).Value = '【Stable Diffusion】参考動画に沿って約100プロンプト構築'
$ws.Cells.Item(28,3).Value = 'システム開発'
$ws.Cells.Item(28,4).Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Cells.Item(28,5).Value = '期限情報なし'
$ws.Cells.Item(28,6).Value = 'https://www.lancers.jp/work/detail/5432055'
$ws.Cells.Item(28,7).Value = 10

# Re-create hyperlinks on the URL column for each data row
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5431917") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5427010") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5427011") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5431609") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5432620") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5432776") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5432819") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5431786") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.lancers.jp/work/detail/5426598") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F11"), "https://www.lancers.jp/work/detail/5432161") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F12"), "https://www.lancers.jp/work/detail/5432465") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F13"), "https://www.lancers.jp/work/detail/5430176") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F14"), "https://www.lancers.jp/work/detail/5431852") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F15"), "https://www.lancers.jp/work/detail/5432661") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F16"), "https://www.lancers.jp/work/detail/5432622") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F17"), "https://www.lancers.jp/work/detail/5432563") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F18"), "https://www.lancers.jp/work/detail/5427007") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F19"), "https://www.lancers.jp/work/detail/5427009") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F20"), "https://www.lancers.jp/work/detail/5432212") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F21"), "https://www.lancers.jp/work/detail/5432621") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F22"), "https://www.lancers.jp/work/detail/5432712") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F23"), "https://www.lancers.jp/work/detail/5432042") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F24"), "https://www.lancers.jp/work/detail/5432323") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F25"), "https://www.lancers.jp/work/detail/5432425") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F26"), "https://www.lancers.jp/work/detail/5432362") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F27"), "https://www.lancers.jp/work/detail/5432305") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F28"), "https://www.lancers.jp/work/detail/5432055") | Out-Null

# Widen the price column to fit the new longer price/period text
$ws.Columns.Item(4).ColumnWidth = 42.14